$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two extra product rows (Sugar row 3 / Chicken row 4).
# Clear() (content + formatting) so the exporter drops the now fully-blank
# rows entirely instead of leaving placeholder <row> elements, while NOT
# shifting the later placeholder rows (21..1000) up.
$ws.Range("A3:E4").Clear()

# Update the remaining product row (row 2): Sugar 5KG -> BEEF 1KG, new qty/price
$ws.Range("A2").Value = "BEEF 1KG"
$ws.Range("B2").Value = 10.0
$ws.Range("C2").Value = 2300.0

# Remark cell (E2) becomes an empty string while staying a text cell.
# A leading apostrophe forces text-typing of an otherwise-empty value.
$ws.Range("E2").Value = "'"

# Match E2's formatting to B2's (drops the extra alignment/quote-prefix
# styling so it lines up with the plain numeric-row style).
$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Narrow columns D and E to match the smaller remark/status content now shown.
# (ColumnWidth is quantized to the nearest 1/6 character by the host, so these
# land on the closest reachable widths to the target 6.71 / 8.14.)
$ws.Range("D1").ColumnWidth = 5.83
$ws.Range("E1").ColumnWidth = 7.33
